# Filter out-of-stock alternates: Only show items with In Stock or Low Stock status as alternatives
#
# The master quantity list was refreshed against current inventory counts so that
# downstream "alternate item" lookups only surface items that are actually In Stock
# or Low Stock (i.e. the Quantity column reflects the latest stock levels). This
# updates the Quantity values for the affected item rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> updated Quantity (column B) value
$updates = [ordered]@{
    90  = -19
    91  = -59
    97  = -50
    101 = 1450
    108 = 1900
    113 = 700
    158 = 3650
    176 = 7200
    213 = 3750
    220 = 4400
    229 = 3750
    287 = 2950
    291 = 1300
    300 = 1600
    302 = 800
    313 = 1900
    318 = 2500
    333 = -1150
    338 = 1850
    354 = 1300
    359 = 50
    411 = 1300
    451 = 0
    489 = 5800
    490 = 3300
    502 = 1750
    503 = 2050
    524 = 1750
    536 = 1500
    539 = 40850
    558 = 1500
    574 = 5500
    594 = 1400
    602 = 1650
    608 = 0
    626 = 2300
    632 = 4150
    635 = 1300
    641 = 5450
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
